$wb = $excel.ActiveWorkbook

# ALC row 15
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(15, 8).Value = 458.44
$ws.Cells.Item(15, 9).Value = 458.44
$ws.Cells.Item(15, 11).Value = 1375.32
$ws.Cells.Item(15, 13).Value = -1206.32

# ALC row 28
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(28, 8).Value = 688.60974
$ws.Cells.Item(28, 9).Value = 624.88464
$ws.Cells.Item(28, 10).Value = 799.06665
$ws.Cells.Item(28, 11).Value = 624.88464
$ws.Cells.Item(28, 12).Value = 799.06665
$ws.Cells.Item(28, 13).Value = -139.88464
$ws.Cells.Item(28, 14).Value = -1769.06665

# ALC row 62
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(62, 8).Value = 4494.8
$ws.Cells.Item(62, 9).Value = 3129.6
$ws.Cells.Item(62, 10).Value = 4722.3335
$ws.Cells.Item(62, 11).Value = 3129.6
$ws.Cells.Item(62, 12).Value = 4722.3335
$ws.Cells.Item(62, 13).Value = -2505.6
$ws.Cells.Item(62, 14).Value = -5970.3335

# ALC row 65
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(65, 8).Value = 4494.8
$ws.Cells.Item(65, 9).Value = 3129.6
$ws.Cells.Item(65, 10).Value = 4722.3335
$ws.Cells.Item(65, 11).Value = 15648
$ws.Cells.Item(65, 12).Value = 23611.6675
$ws.Cells.Item(65, 13).Value = -12528
$ws.Cells.Item(65, 14).Value = -29851.6675

# ALC row 80
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(80, 8).Value = 368.35
$ws.Cells.Item(80, 9).Value = 299.26666
$ws.Cells.Item(80, 10).Value = 575.6
$ws.Cells.Item(80, 11).Value = 897.79998
$ws.Cells.Item(80, 12).Value = 1726.8
$ws.Cells.Item(80, 13).Value = 100.20002
$ws.Cells.Item(80, 14).Value = -3722.8

# ALC row 83
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(83, 8).Value = 368.35
$ws.Cells.Item(83, 9).Value = 299.26666
$ws.Cells.Item(83, 10).Value = 575.6
$ws.Cells.Item(83, 11).Value = 2693.39994
$ws.Cells.Item(83, 12).Value = 5180.400000000001
$ws.Cells.Item(83, 13).Value = 2298.60006
$ws.Cells.Item(83, 14).Value = -15164.4

# ALC row 92
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(92, 8).Value = 611.56525
$ws.Cells.Item(92, 9).Value = 587.7692
$ws.Cells.Item(92, 10).Value = 642.5
$ws.Cells.Item(92, 11).Value = 587.7692
$ws.Cells.Item(92, 12).Value = 642.5
$ws.Cells.Item(92, 13).Value = 660.2308
$ws.Cells.Item(92, 14).Value = -3138.5

# ALC row 103
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(103, 8).Value = 974.8125
$ws.Cells.Item(103, 9).Value = 792.25
$ws.Cells.Item(103, 10).Value = 1522.5
$ws.Cells.Item(103, 11).Value = 2376.75
$ws.Cells.Item(103, 12).Value = 4567.5
$ws.Cells.Item(103, 13).Value = -1790.75
$ws.Cells.Item(103, 14).Value = -5739.5

# ALC row 111
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(111, 8).Value = 2974.1428
$ws.Cells.Item(111, 9).Value = 3397.4
$ws.Cells.Item(111, 10).Value = 1916
$ws.Cells.Item(111, 11).Value = 10192.2
$ws.Cells.Item(111, 12).Value = 5748
$ws.Cells.Item(111, 13).Value = -7125.200000000001
$ws.Cells.Item(111, 14).Value = -11882

# ALC row 137
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(137, 8).Value = 37038604
$ws.Cells.Item(137, 9).Value = 1215.7646
$ws.Cells.Item(137, 10).Value = 100002160
$ws.Cells.Item(137, 11).Value = 3647.2938
$ws.Cells.Item(137, 12).Value = 300006480
$ws.Cells.Item(137, 13).Value = -1097.2938
$ws.Cells.Item(137, 14).Value = -300011580

# ARM row 2
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 892825
$ws.Cells.Item(2, 9).Value = 1386.4286
$ws.Cells.Item(2, 10).Value = 2452842.5
$ws.Cells.Item(2, 11).Value = 1386.4286
$ws.Cells.Item(2, 12).Value = 2452842.5
$ws.Cells.Item(2, 13).Value = -1273.4286
$ws.Cells.Item(2, 14).Value = -2453068.5

# ARM row 61
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(61, 8).Value = 1841.24
$ws.Cells.Item(61, 9).Value = 1870.9131
$ws.Cells.Item(61, 10).Value = 1500
$ws.Cells.Item(61, 11).Value = 1870.9131
$ws.Cells.Item(61, 12).Value = 1500
$ws.Cells.Item(61, 13).Value = -1658.9131
$ws.Cells.Item(61, 14).Value = -1924

# ARM row 116
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(116, 8).Value = 892825
$ws.Cells.Item(116, 9).Value = 1386.4286
$ws.Cells.Item(116, 10).Value = 2452842.5
$ws.Cells.Item(116, 11).Value = 1386.4286
$ws.Cells.Item(116, 12).Value = 2452842.5
$ws.Cells.Item(116, 13).Value = 907.5714
$ws.Cells.Item(116, 14).Value = -2457430.5

# ARM row 122
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(122, 8).Value = 1887.7826
$ws.Cells.Item(122, 9).Value = 1794.9375
$ws.Cells.Item(122, 10).Value = 2100
$ws.Cells.Item(122, 11).Value = 5384.8125
$ws.Cells.Item(122, 12).Value = 6300
$ws.Cells.Item(122, 13).Value = -2934.8125
$ws.Cells.Item(122, 14).Value = -11200

# ARM row 132
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(132, 8).Value = 1147438
$ws.Cells.Item(132, 9).Value = 1439750.4
$ws.Cells.Item(132, 10).Value = 145224.28
$ws.Cells.Item(132, 11).Value = 4319251.199999999
$ws.Cells.Item(132, 12).Value = 435672.84
$ws.Cells.Item(132, 13).Value = -4316721.199999999
$ws.Cells.Item(132, 14).Value = -440732.84

# ARM row 136
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(136, 8).Value = 1841.24
$ws.Cells.Item(136, 9).Value = 1870.9131
$ws.Cells.Item(136, 10).Value = 1500
$ws.Cells.Item(136, 11).Value = 5612.7393
$ws.Cells.Item(136, 12).Value = 4500
$ws.Cells.Item(136, 13).Value = -3062.7393
$ws.Cells.Item(136, 14).Value = -9600

# BSM row 3
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 892825
$ws.Cells.Item(3, 9).Value = 1386.4286
$ws.Cells.Item(3, 10).Value = 2452842.5
$ws.Cells.Item(3, 11).Value = 1386.4286
$ws.Cells.Item(3, 12).Value = 2452842.5
$ws.Cells.Item(3, 13).Value = -1272.4286
$ws.Cells.Item(3, 14).Value = -2453070.5

# BSM row 55
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(55, 8).Value = 46587.5
$ws.Cells.Item(55, 10).Value = 46587.5
$ws.Cells.Item(55, 12).Value = 46587.5
$ws.Cells.Item(55, 14).Value = -47133.5

# BSM row 80
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(80, 8).Value = 186
$ws.Cells.Item(80, 9).Value = 169.71428
$ws.Cells.Item(80, 10).Value = 195.5
$ws.Cells.Item(80, 11).Value = 169.71428
$ws.Cells.Item(80, 12).Value = 195.5
$ws.Cells.Item(80, 13).Value = 828.28572
$ws.Cells.Item(80, 14).Value = -2191.5

# BSM row 83
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(83, 8).Value = 186
$ws.Cells.Item(83, 9).Value = 169.71428
$ws.Cells.Item(83, 10).Value = 195.5
$ws.Cells.Item(83, 11).Value = 848.5714
$ws.Cells.Item(83, 12).Value = 977.5
$ws.Cells.Item(83, 13).Value = 4143.4286
$ws.Cells.Item(83, 14).Value = -10961.5

# BSM row 107
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(107, 8).Value = 456370.44
$ws.Cells.Item(107, 9).Value = 588738.7
$ws.Cells.Item(107, 10).Value = 2536.5715
$ws.Cells.Item(107, 11).Value = 588738.7
$ws.Cells.Item(107, 12).Value = 2536.5715
$ws.Cells.Item(107, 13).Value = -586818.7
$ws.Cells.Item(107, 14).Value = -6376.5715

# CRP row 99
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(99, 8).Value = 1700
$ws.Cells.Item(99, 9).Value = 1580
$ws.Cells.Item(99, 10).Value = 2000
$ws.Cells.Item(99, 11).Value = 1580
$ws.Cells.Item(99, 12).Value = 2000
$ws.Cells.Item(99, 13).Value = -82
$ws.Cells.Item(99, 14).Value = -4996

# CRP row 107
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(107, 8).Value = 1543665.9
$ws.Cells.Item(107, 9).Value = 1984570.4
$ws.Cells.Item(107, 10).Value = 499.66666
$ws.Cells.Item(107, 11).Value = 1984570.4
$ws.Cells.Item(107, 12).Value = 499.66666
$ws.Cells.Item(107, 13).Value = -1982650.4
$ws.Cells.Item(107, 14).Value = -4339.66666

# CRP row 126
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(126, 8).Value = 1700
$ws.Cells.Item(126, 9).Value = 1580
$ws.Cells.Item(126, 10).Value = 2000
$ws.Cells.Item(126, 11).Value = 4740
$ws.Cells.Item(126, 12).Value = 6000
$ws.Cells.Item(126, 13).Value = -2270
$ws.Cells.Item(126, 14).Value = -10940

# CRP row 134
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(134, 8).Value = 3815.8462
$ws.Cells.Item(134, 9).Value = 3483.8333
$ws.Cells.Item(134, 10).Value = 7800
$ws.Cells.Item(134, 11).Value = 10451.4999
$ws.Cells.Item(134, 12).Value = 23400
$ws.Cells.Item(134, 13).Value = -7916.499899999999
$ws.Cells.Item(134, 14).Value = -28470

# CUL row 17
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(17, 8).Value = 647.1429000000001
$ws.Cells.Item(17, 9).Value = 536
$ws.Cells.Item(17, 10).Value = 925
$ws.Cells.Item(17, 11).Value = 1608
$ws.Cells.Item(17, 12).Value = 2775
$ws.Cells.Item(17, 13).Value = -1439
$ws.Cells.Item(17, 14).Value = -3113

# CUL row 122
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(122, 8).Value = 17007696
$ws.Cells.Item(122, 9).Value = 22222766
$ws.Cells.Item(122, 10).Value = 3970023
$ws.Cells.Item(122, 11).Value = 200004894
$ws.Cells.Item(122, 12).Value = 35730207
$ws.Cells.Item(122, 13).Value = -200002444
$ws.Cells.Item(122, 14).Value = -35735107

# CUL row 131
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(131, 8).Value = 3398.244
$ws.Cells.Item(131, 10).Value = 3029.0405
$ws.Cells.Item(131, 12).Value = 9087.121500000001
$ws.Cells.Item(131, 14).Value = -19167.1215

# GSM row 80
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(80, 8).Value = 3085
$ws.Cells.Item(80, 9).Value = 3085
$ws.Cells.Item(80, 10).Value = 0
$ws.Cells.Item(80, 11).Value = 3085
$ws.Cells.Item(80, 12).Value = 0
$ws.Cells.Item(80, 13).Value = -2087
$ws.Cells.Item(80, 14).ClearContents()

# GSM row 83
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(83, 8).Value = 3085
$ws.Cells.Item(83, 9).Value = 3085
$ws.Cells.Item(83, 10).Value = 0
$ws.Cells.Item(83, 11).Value = 15425
$ws.Cells.Item(83, 12).Value = 0
$ws.Cells.Item(83, 13).Value = -10433
$ws.Cells.Item(83, 14).ClearContents()

# GSM row 122
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(122, 8).Value = 25151
$ws.Cells.Item(122, 9).Value = 5351.2
$ws.Cells.Item(122, 10).Value = 44950.8
$ws.Cells.Item(122, 11).Value = 16053.6
$ws.Cells.Item(122, 12).Value = 134852.4
$ws.Cells.Item(122, 13).Value = -13603.6
$ws.Cells.Item(122, 14).Value = -139752.4

# LTW row 46
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(46, 8).Value = 865.26086
$ws.Cells.Item(46, 9).Value = 905.05
$ws.Cells.Item(46, 10).Value = 600
$ws.Cells.Item(46, 11).Value = 905.05
$ws.Cells.Item(46, 12).Value = 600
$ws.Cells.Item(46, 13).Value = -717.05
$ws.Cells.Item(46, 14).Value = -976

# LTW row 100
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(100, 8).Value = 1670.9656
$ws.Cells.Item(100, 9).Value = 1479.875
$ws.Cells.Item(100, 10).Value = 1906.1538
$ws.Cells.Item(100, 11).Value = 1479.875
$ws.Cells.Item(100, 12).Value = 1906.1538
$ws.Cells.Item(100, 13).Value = -938.875
$ws.Cells.Item(100, 14).Value = -2988.1538

# LTW row 122
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(122, 8).Value = 2046.9584
$ws.Cells.Item(122, 9).Value = 1865.1111
$ws.Cells.Item(122, 10).Value = 2156.0667
$ws.Cells.Item(122, 11).Value = 5595.3333
$ws.Cells.Item(122, 12).Value = 6468.2001
$ws.Cells.Item(122, 13).Value = -3145.3333
$ws.Cells.Item(122, 14).Value = -11368.2001

# WVR row 107
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(107, 8).Value = 296.82144
$ws.Cells.Item(107, 9).Value = 287.83334
$ws.Cells.Item(107, 10).Value = 313
$ws.Cells.Item(107, 11).Value = 863.5000200000001
$ws.Cells.Item(107, 12).Value = 939
$ws.Cells.Item(107, 13).Value = 1056.49998
$ws.Cells.Item(107, 14).Value = -4779

# WVR row 112
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(112, 8).Value = 15000
$ws.Cells.Item(112, 10).Value = 15000
$ws.Cells.Item(112, 12).Value = 15000
$ws.Cells.Item(112, 14).Value = -17954
